$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record was inserted above the existing row 87, pushing
# the former rows 87:109 down to 88:110 (dimension grows from R109 to R110).
$ws.Rows.Item(87).Insert()

# The inserted row comes back blank; give the date cell the same number
# format used throughout column D (custom date/time format, style index 2)
# by copying it from the row immediately below (the old row 87, now 88).
$ws.Range("D87").NumberFormat = $ws.Range("D88").NumberFormat

# Columns that are constant for every record in this block.
$ws.Range("A87").Value = 5
$ws.Range("B87").Value = "Macroferia Regional de Talca"
$ws.Range("C87").Value = "Maule"
$ws.Range("E87").Value = 7
$ws.Range("F87").Value = 100112013
$ws.Range("G87").Value = "Alcachofa"
$ws.Range("R87").Value = "Hortaliza"

# New record's data.
$ws.Range("D87").Value = 44855
$ws.Range("H87").Value = "Madrigal"
$ws.Range("I87").Value = "Primera"
$ws.Range("J87").Value = 2000
$ws.Range("K87").Value = 270
$ws.Range("L87").Value = 270
$ws.Range("M87").Value = 270
$ws.Range("N87").Value = "`$/unidad"
$ws.Range("O87").Value = "Región del Maule"
$ws.Range("P87").Value = 270
$ws.Range("Q87").Value = 1
